$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the existing 12 data rows (A1:I12) into the next 12 rows
# (A13:I24), matching the source data/format exactly.
$ws.Range("A1:I12").Copy() | Out-Null
$ws.Range("A13").PasteSpecial() | Out-Null

# Clear clipboard marching ants / mode after pasting.
$excel.CutCopyMode = $false

# Update the active selection to match the new state (K7).
$ws.Range("K7").Select() | Out-Null
